# Refactorizacion archivos .feature, se elimina la linea de cierre de sesion
# The last column ("opcionAutenticacion" / "CONSULTAR_PRODUCTO") represented
# the logout/session-close step of the data-driven test sheet. Remove its
# header and value, clearing the session-close line from the dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$m1 = $ws.Range("M1")
$m2 = $ws.Range("M2")

# Drop the stored values (this also prunes the now-unused shared strings
# "opcionAutenticacion" and "CONSULTAR_PRODUCTO").
$m2.ClearContents()
$m1.ClearContents()

# M1 kept its bold "header" font but loses the yellow highlight fill that
# marked it as the session-close column. Re-derive the formatting from a
# neighboring data-style cell (matching border/fill) and restore the bold
# header font on top of it.
$b2 = $ws.Range("B2")
$b2.Copy()
$m1.PasteSpecial(-4122) # xlPasteFormats
$m1.Font.Bold = $true
$m1.Font.Size = 12
$m1.Interior.Pattern = -4142 # xlNone

# Re-select near the now-trimmed last column.
$ws.Range("E1").Select()
$ws.Range("M1:M2").Select()
